$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.20881881990158
$ws.Range("C2").Value = 0.3264355122434779
$ws.Range("D2").Value = 0.02935279706186122
$ws.Range("F2").Value = 0.344805483763615
$ws.Range("G2").Value = 0.002377528453549117
$ws.Range("O2").Value = 1.079923893146457

$ws.Range("B3").Value = 1.059706637021634
$ws.Range("C3").Value = 0.2878067191863067
$ws.Range("D3").Value = 0.02577375693846307
$ws.Range("F3").Value = 0.3498603769831341
$ws.Range("G3").Value = 0.002380122956698119
$ws.Range("O3").Value = 1.107630825696788

$ws.Range("B4").Value = 0.967833751113119
$ws.Range("C4").Value = 0.2639920444177619
$ws.Range("D4").Value = 0.02356741692385356
$ws.Range("F4").Value = 0.3534422693849635
$ws.Range("G4").Value = 0.002381799108162623
$ws.Range("O4").Value = 1.1262232286211

$ws.Range("B5").Value = 0.9303174010310045
$ws.Range("C5").Value = 0.2542638457374551
$ws.Range("D5").Value = 0.02266616816404365
$ws.Range("F5").Value = 0.3550216891857012
$ws.Range("G5").Value = 0.002382503114638015
$ws.Range("O5").Value = 1.134195512135946

$ws.Range("B6").Value = 0.9240832362210654
$ws.Range("C6").Value = 0.2526470863911641
$ws.Range("D6").Value = 0.02251638870773576
$ws.Range("F6").Value = 0.3552911706230226
$ws.Range("G6").Value = 0.002382621281933414
$ws.Range("O6").Value = 1.135543156739175

$ws.Range("B7").Value = 0.9673281026824156
$ws.Range("C7").Value = 0.2638609408130606
$ws.Range("D7").Value = 0.0235552709761393
$ws.Range("F7").Value = 0.3534630856744769
$ws.Range("G7").Value = 0.002381808517379318
$ws.Range("O7").Value = 1.126329145507725

$ws.Range("B8").Value = 1.157472209034779
$ws.Range("C8").Value = 0.3131367412880763
$ws.Range("D8").Value = 0.02812060668412641
$ws.Range("F8").Value = 0.3464488854397416
$ws.Range("G8").Value = 0.002378405824241845
$ws.Range("O8").Value = 1.089148287501274

$ws.Range("B9").Value = 1.527738573242232
$ws.Range("C9").Value = 0.4089742302664945
$ws.Range("D9").Value = 0.03700105409242838
$ws.Range("F9").Value = 0.3365093123109659
$ws.Range("G9").Value = 0.002372389838385686
$ws.Range("O9").Value = 1.028847679007697

$ws.Range("B10").Value = 1.798095111000748
$ws.Range("C10").Value = 0.4788739395010566
$ws.Range("D10").Value = 0.04347898477602996
$ws.Range("F10").Value = 0.3315619779370209
$ws.Range("G10").Value = 0.002368366278263629
$ws.Range("O10").Value = 0.9923313373582658

$ws.Range("B11").Value = 1.920705992656167
$ws.Range("C11").Value = 0.5105563203329666
$ws.Range("D11").Value = 0.04641537002339646
$ws.Range("F11").Value = 0.3298290132429074
$ws.Range("G11").Value = 0.002366621086280486
$ws.Range("O11").Value = 0.9774308670322114

$ws.Range("B12").Value = 1.967079552613427
$ws.Range("C12").Value = 0.5225364181427494
$ws.Range("D12").Value = 0.04752574486988692
$ws.Range("F12").Value = 0.3292477186641207
$ws.Range("G12").Value = 0.002365972408052262
$ws.Range("O12").Value = 0.9720362194381096

$ws.Range("B13").Value = 1.957094729976006
$ws.Range("C13").Value = 0.5199570723517581
$ws.Range("D13").Value = 0.04728667668278774
$ws.Range("F13").Value = 0.3293695697030259
$ws.Range("G13").Value = 0.002366111571420914
$ws.Range("O13").Value = 0.9731870021451101

$ws.Range("B14").Value = 1.924522320302117
$ws.Range("C14").Value = 0.5115422817900708
$ws.Range("D14").Value = 0.04650675308538155
$ws.Range("F14").Value = 0.3297796849319141
$ws.Range("G14").Value = 0.00236656747502165
$ws.Range("O14").Value = 0.9769820688298694

$ws.Range("B15").Value = 1.904563344613166
$ws.Range("C15").Value = 0.5063857004712418
$ws.Range("D15").Value = 0.0460288208671642
$ws.Range("F15").Value = 0.3300406670653899
$ws.Range("G15").Value = 0.002366848315194093
$ws.Range("O15").Value = 0.9793389857012329

$ws.Range("B16").Value = 1.790074427838363
$ws.Range("C16").Value = 0.476801036386064
$ws.Range("D16").Value = 0.04328686893724409
$ws.Range("F16").Value = 0.331685690484818
$ws.Range("G16").Value = 0.002368482039072286
$ws.Range("O16").Value = 0.9933397016195045

$ws.Range("B17").Value = 1.719741181059476
$ws.Range("C17").Value = 0.4586217266189578
$ws.Range("D17").Value = 0.04160204403396506
$ws.Range("F17").Value = 0.3328277939723421
$ws.Range("G17").Value = 0.002369506044121633
$ws.Range("O17").Value = 1.002368196198589

$ws.Range("B18").Value = 1.679252146082376
$ws.Range("C18").Value = 0.448154655376527
$ws.Range("D18").Value = 0.04063199690197905
$ws.Range("F18").Value = 0.3335333905351945
$ws.Range("G18").Value = 0.002370103042051
$ws.Range("O18").Value = 1.007722087565639

$ws.Range("B19").Value = 1.665537283560241
$ws.Range("C19").Value = 0.4446088506086312
$ws.Range("D19").Value = 0.04030338926853005
$ws.Range("F19").Value = 0.3337806401480208
$ws.Range("G19").Value = 0.002370306554219607
$ws.Range("O19").Value = 1.009562409907417

$ws.Range("B20").Value = 1.727231937727595
$ws.Range("C20").Value = 0.4605580697222535
$ws.Range("D20").Value = 0.04178149847527379
$ws.Range("F20").Value = 0.3327011722829596
$ws.Range("G20").Value = 0.00236939620757547
$ws.Range("O20").Value = 1.001390427794306

$ws.Range("B21").Value = 1.934091180301095
$ws.Range("C21").Value = 0.5140143860404578
$ws.Range("D21").Value = 0.04673587870271945
$ws.Range("F21").Value = 0.3296571861987303
$ws.Range("G21").Value = 0.002366433234545498
$ws.Range("O21").Value = 0.9758606252294157

$ws.Range("B22").Value = 2.068955072803988
$ws.Range("C22").Value = 0.5488498142270259
$ws.Range("D22").Value = 0.04996466797548749
$ws.Range("F22").Value = 0.3281048169948662
$ws.Range("G22").Value = 0.002364567776345589
$ws.Range("O22").Value = 0.9606208371577338

$ws.Range("B23").Value = 1.997006741758582
$ws.Range("C23").Value = 0.530267009348222
$ws.Range("D23").Value = 0.04824226460837622
$ws.Range("F23").Value = 0.328893187346452
$ws.Range("G23").Value = 0.002365556926856003
$ws.Range("O23").Value = 0.9686217512420967

$ws.Range("B24").Value = 1.723845533144129
$ws.Range("C24").Value = 0.4596826973699422
$ws.Range("D24").Value = 0.04170037153639328
$ws.Range("F24").Value = 0.3327582653996402
$ws.Range("G24").Value = 0.002369445839116024
$ws.Range("O24").Value = 1.001831968388871

$ws.Range("B25").Value = 1.427860052413848
$ws.Range("C25").Value = 0.3831355107910781
$ws.Range("D25").Value = 0.03460664892247678
$ws.Range("F25").Value = 0.3387867327414185
$ws.Range("G25").Value = 0.002373947428706409
$ws.Range("O25").Value = 1.043800291419686
